$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.008.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.788.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.555'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.99'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.280'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0659'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0928'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.045.79'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.798.92'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.634'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.006.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '252.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0741'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('E25').Value = '  -2.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.65%  '
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.23%  '
$ws.Range('E29').Value = '  -2.63%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.82'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0515'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.60'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.451.57'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.62%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.32'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.82'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.65%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.899'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0512'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.947.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.82%  '
